# Re-label Lang 67 and Lang 68 as Lang 01 and Lang 02, then re-sort each
# block (the "success" rows and the "fail" rows) by the LANG column so the
# newly-renamed hypo01/hypo02 rows take their place at the top of each
# block in ascending order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Re-label the hypotheses -------------------------------------------------
# "success" block occupies rows 2-67, "fail" block occupies rows 68-133.
# hypo67 -> hypo01 and hypo68 -> hypo02 in both blocks.
$ws.Range("A66").Value = "hypo01"
$ws.Range("A67").Value = "hypo02"
$ws.Range("A132").Value = "hypo01"
$ws.Range("A133").Value = "hypo02"

# --- 2. Re-sort each block ascending by the LANG column (A), keeping the
#        SORT (B) and COUNT (C) columns aligned with their row. -----------------
$successRange = $ws.Range("A2:C67")
$successKey = $ws.Range("A2:A67")
$successRange.Sort($successKey, 1)

$failRange = $ws.Range("A68:C133")
$failKey = $ws.Range("A68:A133")
$failRange.Sort($failKey, 1)

# --- 3. Leave the active selection on G3, matching the saved workbook state. ---
$ws.Range("G3").Select() | Out-Null
